$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ATS Accuracy")
$ws1.Range("B2").Value = 4
$ws1.Range("D2").Value = 63
$ws1.Range("E2").Value = 93.7

$ws1.Range("B3").Value = 3
$ws1.Range("D3").Value = 64
$ws1.Range("E3").Value = 95.3

$ws1.Range("B4").Value = 5
$ws1.Range("D4").Value = 19
$ws1.Range("E4").Value = 73.7

$ws2 = $wb.Worksheets.Item("Total Accuracy")
$ws2.Range("C2").Value = 60
$ws2.Range("D2").Value = 64
$ws2.Range("E2").Value = 93.8

$ws2.Range("B3").Value = 3
$ws2.Range("C3").Value = 55
$ws2.Range("E3").Value = 94.8

$ws2.Range("B5").Value = 3
$ws2.Range("D5").Value = 8
$ws2.Range("E5").Value = 62.5
